$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-23: update Price (D) and Volume(1h) (E) columns only
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.974.58'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.906.35'
$ws.Range('E3').Value = '  -3.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.19'
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4594'
$ws.Range('E7').Value = '  -1.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3824'
$ws.Range('E8').Value = '  -2.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07711'
$ws.Range('E9').Value = '  -2.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9796'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.07'
$ws.Range('E11').Value = '  -2.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.942.74'
$ws.Range('E12').Value = '  -1.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.938'
$ws.Range('E13').Value = '  -3.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.659'
$ws.Range('E14').Value = '  -2.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07031'
$ws.Range('E15').Value = '  -0.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.004'
$ws.Range('E16').Value = '  +0.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '83.75'
$ws.Range('E17').Value = '  -4.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009460'
$ws.Range('E18').Value = '  -4.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.64'
$ws.Range('E19').Value = '  -3.56%  '
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '28.952.11'
$ws.Range('E21').Value = '  -1.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.298'
$ws.Range('E22').Value = '  -4.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.88'
$ws.Range('E23').Value = '  -2.20%  '

# Rows 24-51: coin list shifted; update Coin, Link, Price, Volume(1h)
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.094'
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.08'
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '19.06'
$ws.Range('E26').Value = '  -2.32%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.648'
$ws.Range('E27').Value = '  -2.25%  '
$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '117.29'
$ws.Range('E28').Value = '  -1.88%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.839'
$ws.Range('E29').Value = '  -3.34%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09270'
$ws.Range('E30').Value = '  -1.32%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.8642'
$ws.Range('E31').Value = '  -2.97%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.078'
$ws.Range('E32').Value = '  -2.85%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.247'
$ws.Range('E33').Value = '  -5.25%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.947'
$ws.Range('E34').Value = '  -6.91%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.05711'
$ws.Range('E35').Value = '  -1.74%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.145'
$ws.Range('E36').Value = '  -2.08%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.002'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02041'
$ws.Range('E38').Value = '  -2.73%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5491'
$ws.Range('E39').Value = '  -3.82%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.387'
$ws.Range('E40').Value = '  -4.55%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1755'
$ws.Range('E41').Value = '  -2.16%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.302'
$ws.Range('E42').Value = '  -3.48%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.769'
$ws.Range('E43').Value = '  +0.46%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5168'
$ws.Range('E44').Value = '  -3.16%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.22'
$ws.Range('E45').Value = '  -4.31%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.06822'
$ws.Range('E46').Value = '  -1.30%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.058'
$ws.Range('E47').Value = '  -6.23%  '
$ws.Range('B48').Value = 'PEPE'
$ws.Range('C48').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000002570'
$ws.Range('E48').Value = '  -16.76%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '110.38'
$ws.Range('E49').Value = '  -2.80%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.771'
$ws.Range('E50').Value = '  -3.07%  '
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.001'
$ws.Range('E51').Value = '  -0.03%  '
